$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'8.64%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'18.96%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.338"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'6.33%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08166"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'8.41%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.598"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.26%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.210"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'31.30%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'13.48%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'7.03%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09657"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.67%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'7.89%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.20%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001324"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.76%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005817"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.377"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.08%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.440"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.62%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3394"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.00%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'8.149"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.14%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1419"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.75%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2912"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-9.59%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04301"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.53%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001304"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.02%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004270"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'9.10%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001350"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'9.65%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003539"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-4.97%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02762"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'14.76%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05594"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'7.44%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006299"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.08%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007690"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.66%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1448"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'9.23%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007680"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.91%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008103"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'4.27%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3191"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.48%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006959"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'5.19%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'36.49%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.85%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.08%"
$ws.Range("E51").Style = "Normal"
